$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (rich-text shared strings) ---
$ws.Range("A8").Value = "Volume 29   Number  50"
$ws.Range("C9").Value = "Report Covering the Week  12/12/2022  Through  12/18/2022"

# --- Precinct crime-stat table updates (rows 14-29) ---
$data = @{
  "N14" = -88.095238095238
  "F15" = 4
  "G15" = 3
  "H15" = 33.333333333333
  "I15" = 30
  "J15" = 31
  "K15" = -3.225806451612
  "L15" = -16.666666666666
  "M15" = 50
  "N15" = -6.25
  "C16" = 3
  "D16" = 6
  "E16" = -50
  "F16" = 22
  "G16" = 22
  "H16" = 0
  "I16" = 295
  "J16" = 193
  "K16" = 52.849740932642
  "L16" = 53.645833333333
  "M16" = 6.498194945848
  "N16" = -64.670658682634
  "C17" = 6
  "D17" = 6
  "E17" = 0
  "F17" = 27
  "G17" = 24
  "H17" = 12.5
  "I17" = 413
  "J17" = 305
  "K17" = 35.409836065573
  "L17" = 65.2
  "M17" = 7.272727272727
  "N17" = -16.396761133603
  "C18" = 2
  "D18" = 8
  "E18" = -75
  "F18" = 8
  "G18" = 17
  "H18" = -52.941176470588
  "I18" = 154
  "J18" = 107
  "K18" = 43.925233644859
  "L18" = 5.479452054794
  "M18" = -24.878048780487
  "N18" = -84.857423795476
  "C19" = 5
  "D19" = 7
  "E19" = -28.571428571428
  "F19" = 28
  "H19" = -22.222222222222
  "I19" = 386
  "J19" = 326
  "K19" = 18.40490797546
  "L19" = 33.103448275862
  "M19" = 45.112781954887
  "N19" = -11.872146118721
  "C20" = 4
  "D20" = 10
  "E20" = -60
  "F20" = 23
  "G20" = 20
  "H20" = 15
  "I20" = 218
  "J20" = 142
  "K20" = 53.521126760563
  "L20" = 111.650485436893
  "M20" = 186.842105263158
  "N20" = -48.584905660377
  "C21" = 21
  "D21" = 38
  "E21" = -44.736842105263
  "F21" = 112
  "G21" = 122
  "H21" = -8.196721311475
  "I21" = 1501
  "J21" = 1115
  "K21" = 34.618834080717
  "L21" = 46.868884540117
  "M21" = 21.735604217356
  "N21" = -54.265691651432
  "D22" = 3
  "G22" = 8
  "H22" = -75
  "J22" = 34
  "K22" = -14.705882352941
  "L22" = -3.333333333333
  "M22" = 16
  "D23" = 1
  "G23" = 5
  "J23" = 22
  "K23" = -27.272727272727
  "L23" = 0
  "M23" = 23.076923076923
  "C24" = 14
  "D24" = 11
  "E24" = 27.272727272727
  "F24" = 66
  "G24" = 50
  "H24" = 32
  "I24" = 786
  "J24" = 543
  "K24" = 44.751381215469
  "L24" = 23.390894819466
  "M24" = 20.73732718894
  "C25" = 6
  "D25" = 10
  "E25" = -40
  "F25" = 36
  "G25" = 33
  "H25" = 9.090909090909
  "I25" = 522
  "J25" = 386
  "K25" = 35.233160621761
  "L25" = 47.875354107648
  "M25" = -17.795275590551
  "D26" = 2
  "E26" = -50
  "G26" = 4
  "H26" = 25
  "I26" = 53
  "J26" = 61
  "K26" = -13.11475409836
  "L26" = -22.058823529411
  "C27" = 1
  "D27" = 1
  "E27" = 0
  "F27" = 10
  "G27" = 6
  "H27" = 66.666666666666
  "I27" = 95
  "J27" = 105
  "K27" = -9.523809523809
  "L27" = 72.727272727272
  "G28" = 1
  "L28" = 0
  "G29" = 1
  "L29" = -9.523809523809
}

foreach ($key in $data.Keys) {
  $ws.Range($key).Value = $data[$key]
}

# C18 changes from a text placeholder to a genuine numeric cell;
# restore the numeric (#,##0) format to match the rest of column C.
$ws.Range("C18").NumberFormat = "#,##0"

Write-Host "edit complete"
